# Update the two-digit multiplication problems in the single table of the
# document. Cell text is set positionally (row/column) rather than via a
# global Find/Replace because several of the new values coincide with other
# cells' old values, which would make a naive text-based replace ambiguous
# and order-dependent.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, col) -> new expression text. Rows with problems are 1, 5, 10,
# 15 and 20 (1-based); every row has 5 columns.
$updates = @(
    @{ Row = 1;  Col = 1; Old = "84×26="; New = "17×53=" },
    @{ Row = 1;  Col = 2; Old = "77×71="; New = "77×74=" },
    @{ Row = 1;  Col = 3; Old = "93×20="; New = "70×70=" },
    @{ Row = 1;  Col = 4; Old = "11×20="; New = "59×64=" },
    @{ Row = 1;  Col = 5; Old = "83×26="; New = "84×15=" },

    @{ Row = 5;  Col = 1; Old = "37×70="; New = "49×86=" },
    @{ Row = 5;  Col = 2; Old = "53×66="; New = "95×25=" },
    @{ Row = 5;  Col = 3; Old = "57×13="; New = "18×50=" },
    @{ Row = 5;  Col = 4; Old = "56×98="; New = "64×51=" },
    @{ Row = 5;  Col = 5; Old = "23×22="; New = "67×15=" },

    @{ Row = 10; Col = 1; Old = "65×73="; New = "64×41=" },
    @{ Row = 10; Col = 2; Old = "79×95="; New = "15×40=" },
    @{ Row = 10; Col = 3; Old = "68×86="; New = "95×74=" },
    @{ Row = 10; Col = 4; Old = "24×44="; New = "31×37=" },
    @{ Row = 10; Col = 5; Old = "64×41="; New = "58×18=" },

    @{ Row = 15; Col = 1; Old = "47×76="; New = "36×18=" },
    @{ Row = 15; Col = 2; Old = "94×41="; New = "42×95=" },
    @{ Row = 15; Col = 3; Old = "19×27="; New = "17×25=" },
    @{ Row = 15; Col = 4; Old = "28×46="; New = "49×58=" },
    @{ Row = 15; Col = 5; Old = "83×67="; New = "60×35=" },

    @{ Row = 20; Col = 1; Old = "25×63="; New = "64×28=" },
    @{ Row = 20; Col = 2; Old = "24×98="; New = "74×14=" },
    @{ Row = 20; Col = 3; Old = "78×69="; New = "21×88=" },
    @{ Row = 20; Col = 4; Old = "19×32="; New = "80×27=" },
    @{ Row = 20; Col = 5; Old = "33×62="; New = "80×90=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
